$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '311.42'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.54%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '37.72'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-3.79%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.086'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.82%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07786'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-4.18%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.346'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-3.25%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.900'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-3.03%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '8.211'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.90%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.893'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-8.69%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9146'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.61%'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-8.81%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1914'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-2.61%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09275'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2.96%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03405'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-2.39%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001367'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.87%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005764'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-8.06%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.553'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.50%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.030'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.37%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.53%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '3.85%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.02101'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '5,581.60%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.04356'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.41%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001211'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-2.72%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004258'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-9.88%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02120'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-4.27%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04969'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-4.77%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007679'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.85%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009849'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-5.15%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.55%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002058'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-2.18%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008785'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-3.61%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006683'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.42%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.30%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003038'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '0.74%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002097'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.30%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001997'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.30%'
